# Update cryptocurrency price (column D) and 1h volume/change (column E)
# figures on the active worksheet to match the latest scrape.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '29.370.29'
$ws.Range("D3").Value = '1.881.95'
$ws.Range("E3").Value = '  +0.29%  '
$ws.Range("D4").Value = '''1.001'
$ws.Range("E4").Value = '  +0.16%  '
$ws.Range("D5").Value = '''0.7123'
$ws.Range("E5").Value = '  +0.09%  '
$ws.Range("D6").Value = '''242.97'
$ws.Range("E6").Value = '  +0.33%  '
$ws.Range("D8").Value = '''0.08036'
$ws.Range("E8").Value = '  +3.15%  '
$ws.Range("D9").Value = '''0.3172'
$ws.Range("E9").Value = '  +1.72%  '
$ws.Range("E10").Value = '  -0.28%  '
$ws.Range("D11").Value = '''0.08339'
$ws.Range("E11").Value = '  -1.31%  '
$ws.Range("D12").Value = '1.891.08'
$ws.Range("E12").Value = '  +0.96%  '
$ws.Range("D13").Value = '''5.266'
$ws.Range("E13").Value = '  +0.62%  '
$ws.Range("D14").Value = '''94.92'
$ws.Range("E14").Value = '  +4.09%  '
$ws.Range("D15").Value = '''0.7185'
$ws.Range("D16").Value = '''6.375'
$ws.Range("E16").Value = '  +5.34%  '
$ws.Range("D17").Value = '''0.000008676'
$ws.Range("E17").Value = '  +5.40%  '
$ws.Range("D18").Value = '29.374.94'
$ws.Range("E18").Value = '  -0.04%  '
$ws.Range("D19").Value = '''243.30'
$ws.Range("E19").Value = '  +0.90%  '
$ws.Range("D20").Value = '2.146.11'
$ws.Range("E20").Value = '  +1.41%  '
$ws.Range("E21").Value = '  +0.52%  '
$ws.Range("D22").Value = '''1.001'
$ws.Range("E22").Value = '  +0.16%  '
$ws.Range("D23").Value = '''7.831'
$ws.Range("E23").Value = '  +0.59%  '
$ws.Range("E24").Value = '  +0.15%  '
$ws.Range("D25").Value = '''0.1572'
$ws.Range("E25").Value = '  -1.39%  '
$ws.Range("D26").Value = '''9.105'
$ws.Range("E26").Value = '  +0.34%  '
$ws.Range("E27").Value = '  -0.05%  '
$ws.Range("D28").Value = '''18.60'
$ws.Range("E28").Value = '  +0.34%  '
$ws.Range("D29").Value = '''1.514'
$ws.Range("E29").Value = '  +0.05%  '
$ws.Range("D30").Value = '''4.441'
$ws.Range("E30").Value = '  +0.37%  '
$ws.Range("D31").Value = '''4.354'
$ws.Range("E31").Value = '  +0.54%  '
$ws.Range("E32").Value = '  -6.57%  '
$ws.Range("D33").Value = '''0.05410'
$ws.Range("E33").Value = '  +2.22%  '
$ws.Range("D34").Value = '''1.946'
$ws.Range("E34").Value = '  +0.15%  '
$ws.Range("D35").Value = '''0.7744'
$ws.Range("E35").Value = '  +3.93%  '
$ws.Range("E36").Value = '  +0.66%  '
$ws.Range("D37").Value = '''2.686'
$ws.Range("E37").Value = '  -0.27%  '
$ws.Range("D38").Value = '''0.01889'
$ws.Range("E38").Value = '  +0.91%  '
$ws.Range("D39").Value = '1.271.61'
$ws.Range("E39").Value = '  +3.27%  '
$ws.Range("E40").Value = '  +0.88%  '
$ws.Range("D41").Value = '''6.517'
$ws.Range("E41").Value = '  +0.59%  '
$ws.Range("D42").Value = '''0.9176'
$ws.Range("E42").Value = '  +2.73%  '
$ws.Range("D43").Value = '''113.29'
$ws.Range("E43").Value = '  +2.46%  '
$ws.Range("D44").Value = '''74.70'
$ws.Range("E44").Value = '  +2.76%  '
$ws.Range("E46").Value = '  +5.47%  '
$ws.Range("D47").Value = '2.040.49'
$ws.Range("E47").Value = '  +1.14%  '
$ws.Range("D48").Value = '''1.816'
$ws.Range("E48").Value = '  -0.03%  '
$ws.Range("D49").Value = '''0.5223'
$ws.Range("E49").Value = '  +0.19%  '
$ws.Range("D50").Value = '''9.563'
$ws.Range("E50").Value = '  +1.67%  '
$ws.Range("E51").Value = '  +1.31%  '
